# Table4SI.docx edit:
#  - widen columns 2-7 of the table grid (+~134 twips each)
#  - bump header row height 438 -> 476 twips
#  - shrink the "European hake" row height 434 -> 433 twips
#  - corrupt the degree-sign in the six "(ºC)" header labels to "(ÂºC)"
#    (UTF-8 bytes of º re-decoded as Windows-1252 -> "Â" + "º")

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- column widths (w:gridCol, values are twips; Word COM Width is in points) ---
$t.Columns.Item(2).Width = 2481 / 20
$t.Columns.Item(3).Width = 1836 / 20
$t.Columns.Item(4).Width = 1814 / 20
$t.Columns.Item(5).Width = 1925 / 20
$t.Columns.Item(6).Width = 1636 / 20
$t.Columns.Item(7).Width = 1747 / 20

# --- row heights (w:trHeight, values are twips) ---
$t.Rows.Item(1).Height = 476 / 20
$t.Rows.Item(2).Height = 433 / 20

# --- mojibake the degree sign in the six "(ºC)" cells ---
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("(" + [char]0x00BA + "C)", $true, $false, $false, $false, $false, $true, 1, $false, "(" + [char]0x00C2 + [char]0x00BA + "C)", 2)
